# Add season-record columns (Wins, Losses, Ties) after the existing
# table, mirroring the team's overall W-L-T record onto every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns AD, AE, AF = 30, 31, 32
$colWins   = 30
$colLosses = 31
$colTies   = 32

# Header row: copy the existing header style (from A1, style index 1 -
# bold, centered, bordered) onto the three new header cells, then set
# their text.
$ws.Range("A1").Copy()
$ws.Cells.Item(1, $colWins).PasteSpecial(-4122)
$ws.Cells.Item(1, $colLosses).PasteSpecial(-4122)
$ws.Cells.Item(1, $colTies).PasteSpecial(-4122)

$ws.Cells.Item(1, $colWins).Value = "Wins"
$ws.Cells.Item(1, $colLosses).Value = "Losses"
$ws.Cells.Item(1, $colTies).Value = "Ties"

# Data rows 2..45: the team finished the season 83-79-0, so every
# player row gets that same season record.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, $colWins).Value = 83
    $ws.Cells.Item($row, $colLosses).Value = 79
    $ws.Cells.Item($row, $colTies).Value = 0
}
